# Adds two new worksheets ("Status" and "Type") after the existing "Crit"
# sheet, each holding a small two-row-group summary table (current +
# prior snapshot) that feeds the "Status" and "Type" graphs referenced in
# the commit message. Mirrors the layout/formula pattern already used on
# the "Crit" sheet (A col mirrors the M col via formula, row 3/8 hold
# Sum() totals, rows 1/6 are headers).

$wb = $excel.ActiveWorkbook
$crit = $wb.Worksheets.Item("Crit")

# ---------------------------------------------------------------------
# Status sheet
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Add($null, $crit)
$status.Name = "Status"

$status.Range("A1").Formula = "=M1"
$status.Range("M1").Value = "program"
$status.Range("N1").Value = "No"
$status.Range("O1").Value = "Yes"

$status.Range("A2").Formula = "=M2"
$status.Range("M2").Value = "EXIM"
$status.Range("N2").Value = 4708525908.23
$status.Range("O2").Value = 317549234.54
$status.Range("N2:Q3").NumberFormat = '0.00,,,"B"'
$status.Range("N2:Q3").Font.Name = "Calibri"

$status.Range("A3").Formula = "=M3"
$status.Range("M3").Value = "Grand Total"
$status.Range("N3").Formula = "=Sum(N2:N2)"
$status.Range("O3").Formula = "=Sum(O2:O2)"

$status.Range("A6").Formula = "=M6"
$status.Range("M6").Value = "program"
$status.Range("N6").Value = "No"
$status.Range("O6").Value = "Yes"

$status.Range("A7").Formula = "=M7"
$status.Range("M7").Value = "EXIM"
$status.Range("N7").Value = 5793327846.85168
$status.Range("O7").Value = 394649584.014714
$status.Range("N7:Q8").NumberFormat = '0.00,,,"B"'
$status.Range("N7:Q8").Font.Name = "Calibri"

$status.Range("A8").Formula = "=M8"

$status.Range("B2").Select()
$status.Application.ActiveWindow.Zoom = 100
$status.Application.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# Type sheet
# ---------------------------------------------------------------------
$type = $wb.Worksheets.Add($null, $status)
$type.Name = "Type"

$type.Range("A1").Formula = "=M1"
$type.Range("M1").Value = "program"
$type.Range("N1").Value = "Guarantee"
$type.Range("O1").Value = "Insurance"
$type.Range("P1").Value = "Loan"
$type.Range("Q1").Value = "Working Capital"

$type.Range("A2").Formula = "=M2"
$type.Range("M2").Value = "EXIM"
$type.Range("N2").Value = 850558403.33
$type.Range("O2").Value = 1568098368.04
$type.Range("P2").Value = 1286773861
$type.Range("Q2").Value = 1320644510.4
$type.Range("N2:S3").NumberFormat = '0.00,,,"B"'
$type.Range("N2:S3").Font.Name = "Calibri"

$type.Range("A3").Formula = "=M3"
$type.Range("M3").Value = "Grand Total"
$type.Range("N3").Formula = "=Sum(N2:N2)"
$type.Range("O3").Formula = "=Sum(O2:O2)"
$type.Range("P3").Formula = "=Sum(P2:P2)"
$type.Range("Q3").Formula = "=Sum(Q2:Q2)"

$type.Range("A6").Formula = "=M6"
$type.Range("M6").Value = "program"
$type.Range("N6").Value = "Guarantee"
$type.Range("O6").Value = "Insurance"
$type.Range("P6").Value = "Loan"
$type.Range("Q6").Value = "Working Capital"

$type.Range("A7").Formula = "=M7"
$type.Range("M7").Value = "EXIM"
$type.Range("N7").Value = 1069332501.26643
$type.Range("O7").Value = 1925461736.14809
$type.Range("P7").Value = 1575047154.93259
$type.Range("Q7").Value = 1618136038.51928
$type.Range("N7:S8").NumberFormat = '0.00,,,"B"'
$type.Range("N7:S8").Font.Name = "Calibri"

$type.Range("A8").Formula = "=M8"

$type.Range("B2").Select()
$type.Application.ActiveWindow.Zoom = 100
$type.Application.ActiveWindow.FreezePanes = $true

# Leave the workbook's original active sheet selected, as in the source file.
$crit.Activate()
